# Refresh the cryptos price/volume table (columns D and E) with the latest
# scraped values. Numeric-looking Price strings (column D) are written with
# a leading apostrophe to force Excel to keep them as text (matching the
# original inlineStr storage, e.g. "523.71" must stay text, not become the
# number 523.71), then the cell's original style is restored so we don't
# leave a stray "text" number-format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style_D2 = $ws.Range("D2").Style
$ws.Range("D2").Value = "'58.951.51"
$ws.Range("D2").Style = $style_D2
$ws.Range("E2").Value = '  -2.39%  '
$style_D3 = $ws.Range("D3").Style
$ws.Range("D3").Value = "'2.656.40"
$ws.Range("D3").Style = $style_D3
$ws.Range("E3").Value = '  -1.17%  '
$ws.Range("E4").Value = '  -0.23%  '
$style_D5 = $ws.Range("D5").Style
$ws.Range("D5").Value = "'523.71"
$ws.Range("D5").Style = $style_D5
$ws.Range("E5").Value = '  +0.36%  '
$style_D6 = $ws.Range("D6").Style
$ws.Range("D6").Value = "'144.49"
$ws.Range("D6").Style = $style_D6
$ws.Range("E6").Value = '  -1.16%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("E8").Value = '  -1.23%  '
$ws.Range("E9").Value = '  +8.65%  '
$style_D10 = $ws.Range("D10").Style
$ws.Range("D10").Value = "'0.103"
$ws.Range("D10").Style = $style_D10
$ws.Range("E10").Value = '  -2.69%  '
$ws.Range("E11").Value = '  -1.76%  '
$ws.Range("E12").Value = '  +1.29%  '
$style_D13 = $ws.Range("D13").Style
$ws.Range("D13").Value = "'3.120.66"
$ws.Range("D13").Style = $style_D13
$ws.Range("E13").Value = '  -1.04%  '
$style_D14 = $ws.Range("D14").Style
$ws.Range("D14").Value = "'58.958.37"
$ws.Range("D14").Style = $style_D14
$ws.Range("E14").Value = '  -2.40%  '
$style_D15 = $ws.Range("D15").Style
$ws.Range("D15").Value = "'21.06"
$ws.Range("D15").Style = $style_D15
$ws.Range("E15").Value = '  -1.16%  '
$ws.Range("E16").Value = '  -1.80%  '
$style_D17 = $ws.Range("D17").Style
$ws.Range("D17").Value = "'2.667.25"
$ws.Range("D17").Style = $style_D17
$ws.Range("E17").Value = '  -2.91%  '
$style_D18 = $ws.Range("D18").Style
$ws.Range("D18").Value = "'338.75"
$ws.Range("D18").Style = $style_D18
$ws.Range("E18").Value = '  -3.47%  '
$ws.Range("E19").Value = '  -3.97%  '
$ws.Range("E20").Value = '  -1.72%  '
$ws.Range("E21").Value = '  +0.43%  '
$ws.Range("E22").Value = '  -0.09%  '
$style_D23 = $ws.Range("D23").Style
$ws.Range("D23").Value = "'63.81"
$ws.Range("D23").Style = $style_D23
$ws.Range("E23").Value = '  +0.95%  '
$ws.Range("E24").Value = '  -1.21%  '
$style_D25 = $ws.Range("D25").Style
$ws.Range("D25").Value = "'0.165"
$ws.Range("D25").Style = $style_D25
$ws.Range("E25").Value = '  -1.67%  '
$style_D26 = $ws.Range("D26").Style
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = $style_D26
$ws.Range("E26").Value = '  +0.55%  '
$ws.Range("E27").Value = '  -1.70%  '
$ws.Range("E28").Value = '  -2.63%  '
$style_D29 = $ws.Range("D29").Style
$ws.Range("D29").Value = "'6.66"
$ws.Range("D29").Style = $style_D29
$ws.Range("E29").Value = '  -3.46%  '
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("E31").Value = '  -0.16%  '
$style_D32 = $ws.Range("D32").Style
$ws.Range("D32").Value = "'18.83"
$ws.Range("D32").Style = $style_D32
$ws.Range("E32").Value = '  -1.31%  '
$style_D33 = $ws.Range("D33").Style
$ws.Range("D33").Value = "'149.86"
$ws.Range("D33").Style = $style_D33
$ws.Range("E33").Value = '  +0.79%  '
$ws.Range("E34").Value = '  -4.88%  '
$ws.Range("E35").Value = '  -3.19%  '
$style_D36 = $ws.Range("D36").Style
$ws.Range("D36").Value = "'0.891"
$ws.Range("D36").Style = $style_D36
$ws.Range("E36").Value = '  -6.31%  '
$ws.Range("E37").Value = '  -0.51%  '
$style_D38 = $ws.Range("D38").Style
$ws.Range("D38").Value = "'36.81"
$ws.Range("D38").Style = $style_D38
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("E39").Value = '  -6.44%  '
$ws.Range("E40").Value = '  -2.97%  '
$style_D41 = $ws.Range("D41").Style
$ws.Range("D41").Value = "'0.617"
$ws.Range("D41").Style = $style_D41
$ws.Range("E42").Value = '  +0.36%  '
$ws.Range("E43").Value = '  -0.85%  '
$style_D44 = $ws.Range("D44").Style
$ws.Range("D44").Value = "'275.47"
$ws.Range("D44").Style = $style_D44
$ws.Range("E44").Value = '  -2.62%  '
$ws.Range("E45").Value = '  -2.18%  '
$ws.Range("E46").Value = '  +2.03%  '
$style_D47 = $ws.Range("D47").Style
$ws.Range("D47").Value = "'2.047.12"
$ws.Range("D47").Style = $style_D47
$ws.Range("E47").Value = '  -3.69%  '
$style_D48 = $ws.Range("D48").Style
$ws.Range("D48").Value = "'0.0531"
$ws.Range("D48").Style = $style_D48
$ws.Range("E48").Value = '  -1.88%  '
$ws.Range("E49").Value = '  -3.03%  '
$style_D50 = $ws.Range("D50").Style
$ws.Range("D50").Value = "'18.93"
$ws.Range("D50").Style = $style_D50
$ws.Range("E50").Value = '  -0.86%  '
$ws.Range("E51").Value = '  -2.97%  '
